# SectorGroup.xlsx column re-order.
#
# The codeforiati "group"/"category" columns (D, E, F, G) get rearranged
# so the column order becomes:
#   D = codeforiati:category-name   (was E)
#   E = codeforiati:group-name      (was G)
#   F = codeforiati:category-code   (unchanged)
#   G = codeforiati:group-code      (was D)
#
# i.e. a 3-way rotation of columns D -> G -> E -> D, leaving F untouched.
# This applies uniformly to every row, including the header row.
#
# We use Range.Copy (a true COM cell copy) rather than reading/writing
# .Value, because round-tripping numeric-looking text (e.g. "110") through
# .Value would coerce it to a real number and lose the shared-string/text
# representation. Copy() preserves the original cell type exactly.
#
# A scratch column (far to the right of the used range) is used to stage
# the old column-D contents during the 3-way rotation, then cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count()

$colD = 4
$colE = 5
$colF = 6
$colG = 7
$scratchCol = 9   # column I: two columns past G, outside the used range

for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $colD).Copy($ws.Cells.Item($r, $scratchCol))
}
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $colE).Copy($ws.Cells.Item($r, $colD))
}
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $colG).Copy($ws.Cells.Item($r, $colE))
}
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $scratchCol).Copy($ws.Cells.Item($r, $colG))
}
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $scratchCol).Clear()
}
